# Update "NV-15 Truong Lam Khanh 7-2024" workbook:
#  1. Insert a new sheet "Don phu phau 1" between "Don sale chinh" and "Luong",
#     populated with a new service-order table (phu phau 1 case for NGUYEN... etc.)
#  2. On "Don sale chinh", change G2 ("Nhom dich vu") from "Cac ngoai khoa khac" to "Tiem"
#  3. On "Luong", insert a new "Ung luong" row per co so + two new "Tong luong" style
#     rows, and refresh the numeric results for the period.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Don sale chinh -> fix G2
# ---------------------------------------------------------------------------
$wsMain = $wb.Worksheets.Item(1)
$wsMain.Range("G2").Value = "Tiêm"

# ---------------------------------------------------------------------------
# 2) Insert new sheet "Đơn phụ phẫu 1" right after "Đơn sale chính"
# ---------------------------------------------------------------------------
$wsSub = $wb.Worksheets.Add($null, $wsMain)
$wsSub.Name = "Đơn phụ phẫu 1"

$subHeaders = @(
    "Tiền tố", "Mã dịch vụ", "Ngày thực hiện", "Cơ sở", "Khách hàng",
    "Nguồn khách", "Nhóm dịch vụ", "Tên dịch vụ", "Sale chính", "Đơn giá gốc",
    "Sale phụ", "Upsale", "Đơn giá", "Thanh toán lần đầu", "Trả sau",
    "Đã thanh toán", "Dư nợ", "Bác sĩ 1", "Bác sĩ 2", "Phụ phẫu 1",
    "Phụ phẫu 2", "Công phụ phẫu 1", "Công phụ phẫu 2",
    "Tỉ lệ chiết khấu sale chính", "Tỉ lệ chiết khấu sale phụ",
    "Chiết khấu sale chính", "Chiết khấu sale phụ"
)
for ($i = 0; $i -lt $subHeaders.Length; $i++) {
    $wsSub.Cells.Item(1, $i + 1).Value = $subHeaders[$i]
}

$wsSub.Range("A2").Value = "HD-LUXURY"
$wsSub.Range("B2").Value = 559
# Force text storage so the "MM-DD-YYYY" string is not auto-converted to a
# date serial number (matches the source file's inlineStr type for C2).
$wsSub.Range("C2").NumberFormat = "@"
$wsSub.Range("C2").Value = "07-16-2024"
$wsSub.Range("D2").Value = "CẦN THƠ"
$wsSub.Range("E2").Value = "Lê Thị Kiều"
$wsSub.Range("F2").Value = "Cá nhân"
$wsSub.Range("G2").Value = "Tiểu phẫu"
$wsSub.Range("H2").Value = "Mở góc mắt"
$wsSub.Range("I2").Value = "Lê Đình Hậu"
$wsSub.Range("J2").Value = 8000000
$wsSub.Range("K2").Value = $null
$wsSub.Range("L2").Value = $null
$wsSub.Range("M2").Value = 8000000
$wsSub.Range("N2").Value = 0
$wsSub.Range("O2").Value = 0
$wsSub.Range("P2").Value = 0
$wsSub.Range("Q2").Value = 8000000
$wsSub.Range("R2").Value = "Nguyễn Hoàng Yến Quyên"
$wsSub.Range("S2").Value = $null
$wsSub.Range("T2").Value = "Trương Lâm Khanh"
$wsSub.Range("U2").Value = $null
$wsSub.Range("V2").Value = 50000
$wsSub.Range("W2").Value = 0
$wsSub.Range("X2").Value = 0.13
$wsSub.Range("Y2").Value = 0
$wsSub.Range("Z2").Value = 0
$wsSub.Range("AA2").Value = 0

$wsSub.Range("A3").Value = "Tổng"
$wsSub.Range("B3").Value = 1
$wsSub.Range("C3").Value = $null
$wsSub.Range("D3").Value = $null
$wsSub.Range("E3").Value = $null
$wsSub.Range("F3").Value = $null
$wsSub.Range("G3").Value = $null
$wsSub.Range("H3").Value = $null
$wsSub.Range("I3").Value = $null
$wsSub.Range("J3").Value = 8000000
$wsSub.Range("K3").Value = $null
$wsSub.Range("L3").Value = 0
$wsSub.Range("M3").Value = 8000000
$wsSub.Range("N3").Value = 0
$wsSub.Range("O3").Value = 0
$wsSub.Range("P3").Value = 0
$wsSub.Range("Q3").Value = 8000000
$wsSub.Range("R3").Value = $null
$wsSub.Range("S3").Value = $null
$wsSub.Range("T3").Value = $null
$wsSub.Range("U3").Value = $null
$wsSub.Range("V3").Value = 50000
$wsSub.Range("W3").Value = 0
$wsSub.Range("X3").Value = 0.13
$wsSub.Range("Y3").Value = 0
$wsSub.Range("Z3").Value = 0
$wsSub.Range("AA3").Value = 0

# ---------------------------------------------------------------------------
# 3) "Lương" sheet: refresh figures + add "Ứng lương" / "Tổng lương" rows
# ---------------------------------------------------------------------------
$wsLuong = $wb.Worksheets.Item("Lương")

$wsLuong.Range("B2").Value = 16
$wsLuong.Range("B3").Value = 560000
$wsLuong.Range("B4").Value = 1714285.714285715
# B5 (Chiết khấu sale chính tại CẦN THƠ) / B6 / B7 / B8 stay the same
$wsLuong.Range("B9").Value = 50000
# B10 (Công phụ phẫu 2 tại CẦN THƠ) stays the same

# Insert a brand-new row right after the "CẦN THƠ" block (old row 11, which
# held "Lương cơ bản tại LONG XUYÊN", is pushed down to row 12), so a new
# "Ứng lương tại CẦN THƠ" row can live at row 11.
$wsLuong.Range("A11:B11").Insert()
$wsLuong.Range("A11").Value = "Ứng lương tại CẦN THƠ"
$wsLuong.Range("B11").Value = -3850000

# The LONG XUYÊN block now occupies rows 12-18. Refresh its base-salary
# figure, then insert a new row right after it (old row 18, "Lương cơ bản
# tại SÓC TRĂNG", is pushed down to row 20) for "Ứng lương tại LONG XUYÊN".
$wsLuong.Range("B12").Value = 1142857.142857143
$wsLuong.Range("A19:B19").Insert()
$wsLuong.Range("A19").Value = "Ứng lương tại LONG XUYÊN"
$wsLuong.Range("B19").Value = -0

# The SÓC TRĂNG block now occupies rows 20-26. Refresh its base-salary
# figure, then insert a new row right after it for "Ứng lương tại SÓC
# TRĂNG".
$wsLuong.Range("B20").Value = 1714285.714285715
$wsLuong.Range("A27:B27").Insert()
$wsLuong.Range("A27").Value = "Ứng lương tại SÓC TRĂNG"
$wsLuong.Range("B27").Value = -0

# Finally, append the four grand-total rows.
$wsLuong.Range("A28").Value = "Tổng lương tại CẦN THƠ"
$wsLuong.Range("B28").Value = 2934285.714285715

$wsLuong.Range("A29").Value = "Tổng lương tại LONG XUYÊN"
$wsLuong.Range("B29").Value = 1142857.142857143

$wsLuong.Range("A30").Value = "Tổng lương tại SÓC TRĂNG"
$wsLuong.Range("B30").Value = 1714285.714285715

$wsLuong.Range("A31").Value = "Tổng lương"
$wsLuong.Range("B31").Value = 5791428.571428572
